$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows 14-17: labels + aggregate formulas ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

# --- New average row for column J (k value check) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Copy B14's formatting (bold, size 12, vertically centered) onto B15:B17
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Match the workbook's final on-screen selection
$ws.Range("B17").Select()
